# BAU Min Req EV Sales Perc.xlsx — apply commit:
# "updated deeper decarb LDV EV adoption percentage and updated BPMRESP to
#  set BAU 100 percent by 2035"
#
# Core change: insert a new row 3 ("Post ACC2") in the 'ZEV Waiver States'
# sheet above the existing lookup row (which becomes row 4), and add a new
# row 5 ("Straightline to 2035 (ACC 2)") holding the straight-line ramp to
# 100% EV sales by 2035. Downstream sheets that reference this table are
# fixed up to point at the right rows again.

$wb = $excel.ActiveWorkbook
$zev = $wb.Worksheets.Item("ZEV Waiver States")

# --- 1. Insert a new row above the old row 3 --------------------------------
# This shifts: old row3->4, old row5->6, old row6->7, old row8->9, old row9->10, ...
# Excel auto-updates same-workbook formula references (including cross-sheet
# references), matching the diff's row4/row12/etc. reference shifts.
$zev.Rows.Item(3).Insert()

# --- 2. Populate new row 5 label first: "Straightline to 2035 (ACC 2)" ------
# (establishes shared-string order to match the authored workbook)
$zev.Range("A5").Value = "Straightline to 2035 (ACC 2)"

# --- 2b. Populate new row 3: "Post ACC2" --------------------------------------
$zev.Range("A3").Value = "Post ACC2"

# D3:G3 mirror the (now shifted-down) historical ZEV data row (row 6)
$zev.Range("D3").Formula = "=D6"
$zev.Range("E3").Formula = "=E6"
$zev.Range("F3").Formula = "=F6"
$zev.Range("G3").Formula = "=G6"

# H3:U3 pick up the new straight-line-to-2035 row (row 5)
$zev.Range("H3").Formula = "=H5"
$zev.Range("I3").Formula = "=I5"
$zev.Range("J3").Formula = "=J5"
$zev.Range("K3").Formula = "=K5"
$zev.Range("L3").Formula = "=L5"
$zev.Range("M3").Formula = "=M5"
$zev.Range("N3").Formula = "=N5"
$zev.Range("O3").Formula = "=O5"
$zev.Range("P3").Formula = "=P5"
$zev.Range("Q3").Formula = "=Q5"
$zev.Range("R3").Formula = "=R5"
$zev.Range("S3").Formula = "=S5"
$zev.Range("T3").Formula = "=T5"
$zev.Range("U3").Formula = "=U5"

# V3:AJ3 hold flat at 100% (each column repeats the previous one)
$zev.Range("V3").Formula = "=U3"
$zev.Range("W3").Formula = "=V3"
$zev.Range("X3").Formula = "=W3"
$zev.Range("Y3").Formula = "=X3"
$zev.Range("Z3").Formula = "=Y3"
$zev.Range("AA3").Formula = "=Z3"
$zev.Range("AB3").Formula = "=AA3"
$zev.Range("AC3").Formula = "=AB3"
$zev.Range("AD3").Formula = "=AC3"
$zev.Range("AE3").Formula = "=AD3"
$zev.Range("AF3").Formula = "=AE3"
$zev.Range("AG3").Formula = "=AF3"
$zev.Range("AH3").Formula = "=AG3"
$zev.Range("AI3").Formula = "=AH3"
$zev.Range("AJ3").Formula = "=AI3"

# Style row 3 like the neighbouring highlighted rows, with a high-precision
# number format on the cells pulled from the legacy ZEV data (D3:G3).
$zev.Range("A3:AJ3").Interior.ColorIndex = $zev.Range("A6").EntireRow.Interior.ColorIndex
$zev.Range("B3:AJ3").NumberFormat = "General"
$zev.Range("D3:G3").NumberFormat = "0.0000000000"

Write-Host "Row 3 (Post ACC2) populated"

# --- 3. Populate new row 5: "Straightline to 2035 (ACC 2)" ------------------
$zev.Range("A5").Value = "Straightline to 2035 (ACC 2)"
$zev.Range("A5").Interior.ColorIndex = -4142
$zev.Range("A5").Font.Bold = $false

$zev.Range("H5").Value = 0.118507
$zev.Range("I5").Value = 0.186314
$zev.Range("J5").Value = 0.254121
$zev.Range("K5").Value = 0.321928
$zev.Range("L5").Value = 0.389736
$zev.Range("M5").Value = 0.457543
$zev.Range("N5").Value = 0.52535
$zev.Range("O5").Value = 0.593157
$zev.Range("P5").Value = 0.660964
$zev.Range("Q5").Value = 0.728771
$zev.Range("R5").Value = 0.796579
$zev.Range("S5").Value = 0.864386
$zev.Range("T5").Value = 0.932193
$zev.Range("U5").Value = 1
$zev.Range("V5").Value = 1

Write-Host "Row 5 (Straightline to 2035) populated"

# --- 4. Row 4 (shifted-down lookup row): D4 is intentionally left blank -----
# (the "Post ACC2"/straightline rows now supply that column instead)
$zev.Range("D4").ClearContents()

Write-Host "Row 4 D cleared"

# --- 5. BMRESP-passenger: re-point E2:AJ2 at the new "Post ACC2" row 3 ------
# B2:D2 legitimately track the shifted-down lookup row (now row 4); E2:AJ2
# should keep tracking row 3, which now carries the updated post-ACC2 figures.
$psg = $wb.Worksheets.Item("BMRESP-passenger")
$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")
foreach ($c in $cols) {
    $psg.Range($c + "2").Formula = "='ZEV Waiver States'!" + $c + "3"
}

Write-Host "BMRESP-passenger row 2 re-pointed"
